$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records (date 2021-11-22, serial 44522) were added to the
# weekly Fruta/Hortalizas dataset. They land at row 140:141, so insert two
# blank rows there (this pushes the former rows 140-168 down to 142-170,
# preserving all of their data/formatting).
$ws.Rows("140:141").Insert()

# New row 140: "Primera" quality record for the Melipilla province
$ws.Range("A140").Value2 = 3
$ws.Range("B140").Value2 = "Femacal de La Calera"
$ws.Range("C140").Value2 = "Coquimbo"
$ws.Range("D140").Value2 = 44522
$ws.Range("E140").Value2 = 5
$ws.Range("F140").Value2 = "Fruta"
$ws.Range("G140").Value2 = 100101
$ws.Range("H140").Value2 = "Berries"
$ws.Range("I140").Value2 = 100112025
$ws.Range("J140").Value2 = "Frutilla"
$ws.Range("K140").Value2 = "Sin especificar"
$ws.Range("L140").Value2 = "Primera"
$ws.Range("M140").Value2 = 135
$ws.Range("N140").Value2 = 6000
$ws.Range("O140").Value2 = 6500
$ws.Range("P140").Value2 = 6259
$ws.Range("Q140").Value2 = "`$/bandeja 7 kilos"
$ws.Range("R140").Value2 = "Provincia de Melipilla"
$ws.Range("S140").Value2 = 894
$ws.Range("T140").Value2 = 7

# New row 141: "Segunda" quality record for the Melipilla province
$ws.Range("A141").Value2 = 3
$ws.Range("B141").Value2 = "Femacal de La Calera"
$ws.Range("C141").Value2 = "Coquimbo"
$ws.Range("D141").Value2 = 44522
$ws.Range("E141").Value2 = 5
$ws.Range("F141").Value2 = "Fruta"
$ws.Range("G141").Value2 = 100101
$ws.Range("H141").Value2 = "Berries"
$ws.Range("I141").Value2 = 100112025
$ws.Range("J141").Value2 = "Frutilla"
$ws.Range("K141").Value2 = "Sin especificar"
$ws.Range("L141").Value2 = "Segunda"
$ws.Range("M141").Value2 = 80
$ws.Range("N141").Value2 = 4500
$ws.Range("O141").Value2 = 4500
$ws.Range("P141").Value2 = 4500
$ws.Range("Q141").Value2 = "`$/bandeja 7 kilos"
$ws.Range("R141").Value2 = "Provincia de Melipilla"
$ws.Range("S141").Value2 = 643
$ws.Range("T141").Value2 = 7
